# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 174.09091
$ws.Range("I5").Value = 101.666664
$ws.Range("K5").Value = 101.666664
$ws.Range("M5").Value = 13.333336

$ws.Range("H12").Value = 600.3333
$ws.Range("I12").Value = 1001
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 1001
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = -831
$ws.Range("N12").Value = -740

$ws.Range("H33").Value = 1515289.4
$ws.Range("I33").Value = 1818306.2
$ws.Range("J33").Value = 205.4
$ws.Range("K33").Value = 1818306.2
$ws.Range("L33").Value = 205.4
$ws.Range("M33").Value = -1818077.2
$ws.Range("N33").Value = -663.4

$ws.Range("H112").Value = 1299.5074
$ws.Range("J112").Value = 1299.5074
$ws.Range("L112").Value = 3898.5222
$ws.Range("N112").Value = -6114.522199999999

$ws.Range("H134").Value = 60681.43
$ws.Range("J134").Value = 60681.43
$ws.Range("L134").Value = 60681.43
$ws.Range("N134").Value = -70821.42999999999

$ws.Range("H137").Value = 1194588.6
$ws.Range("I137").Value = 3406760
$ws.Range("J137").Value = 3419.3462
$ws.Range("K137").Value = 10220280
$ws.Range("L137").Value = 10258.0386
$ws.Range("M137").Value = -10217730
$ws.Range("N137").Value = -15358.0386

$ws.Range("H138").Value = 2302.5
$ws.Range("J138").Value = 3246.2778
$ws.Range("L138").Value = 9738.8334
$ws.Range("N138").Value = -20018.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 314367.72
$ws.Range("I74").Value = 507787.2
$ws.Range("J74").Value = 1920.8462
$ws.Range("K74").Value = 507787.2
$ws.Range("L74").Value = 1920.8462
$ws.Range("M74").Value = -506913.2
$ws.Range("N74").Value = -3668.8462

$ws.Range("H77").Value = 314367.72
$ws.Range("I77").Value = 507787.2
$ws.Range("J77").Value = 1920.8462
$ws.Range("K77").Value = 2538936
$ws.Range("L77").Value = 9604.231
$ws.Range("M77").Value = -2534568
$ws.Range("N77").Value = -18340.231

$ws.Range("H80").Value = 35458.668
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 35458.668
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H122").Value = 2567.4736
$ws.Range("I122").Value = 1552.4615
$ws.Range("J122").Value = 4766.6665
$ws.Range("K122").Value = 4657.3845
$ws.Range("L122").Value = 14299.9995
$ws.Range("M122").Value = -2207.3845
$ws.Range("N122").Value = -19199.9995

$ws.Range("H134").Value = 48890
$ws.Range("J134").Value = 48890
$ws.Range("L134").Value = 48890
$ws.Range("N134").Value = -59030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3299.375
$ws.Range("I134").Value = 1467.3334
$ws.Range("K134").Value = 4402.0002
$ws.Range("M134").Value = -1867.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6538071.5
$ws.Range("I16").Value = 12347370
$ws.Range("J16").Value = 2609.75
$ws.Range("K16").Value = 12347370
$ws.Range("L16").Value = 2609.75
$ws.Range("M16").Value = -12347083
$ws.Range("N16").Value = -3183.75

$ws.Range("H31").Value = 192295.62
$ws.Range("I31").Value = 501321.8
$ws.Range("J31").Value = 2665.9092
$ws.Range("K31").Value = 501321.8
$ws.Range("L31").Value = 2665.9092
$ws.Range("M31").Value = -501026.8
$ws.Range("N31").Value = -3255.9092

$ws.Range("H34").Value = 192295.62
$ws.Range("I34").Value = 501321.8
$ws.Range("J34").Value = 2665.9092
$ws.Range("K34").Value = 501321.8
$ws.Range("L34").Value = 2665.9092
$ws.Range("M34").Value = -501119.8
$ws.Range("N34").Value = -3069.9092

$ws.Range("H41").Value = 37589.715
$ws.Range("J41").Value = 37589.715
$ws.Range("L41").Value = 37589.715
$ws.Range("N41").Value = -38445.715

$ws.Range("H58").Value = 2334.2327
$ws.Range("J58").Value = 3650.5557
$ws.Range("L58").Value = 3650.5557
$ws.Range("N58").Value = -4056.5557

$ws.Range("H113").Value = 6538071.5
$ws.Range("I113").Value = 12347370
$ws.Range("J113").Value = 2609.75
$ws.Range("K113").Value = 12347370
$ws.Range("L113").Value = 2609.75
$ws.Range("M113").Value = -12345200
$ws.Range("N113").Value = -6949.75

$ws.Range("H136").Value = 2334.2327
$ws.Range("J136").Value = 3650.5557
$ws.Range("L136").Value = 10951.6671
$ws.Range("N136").Value = -16051.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 12526429
$ws.Range("I107").Value = 456.92307
$ws.Range("K107").Value = 1370.76921
$ws.Range("M107").Value = 549.2307900000001

$ws.Range("H113").Value = 449
$ws.Range("I113").Value = 463.8
$ws.Range("J113").Value = 433.58334
$ws.Range("K113").Value = 1391.4
$ws.Range("L113").Value = 1300.75002
$ws.Range("M113").Value = 778.5999999999999
$ws.Range("N113").Value = -5640.750019999999

$ws.Range("H131").Value = 783.12
$ws.Range("J131").Value = 803.95746
$ws.Range("L131").Value = 2411.87238
$ws.Range("N131").Value = -12491.87238

$ws.Range("H140").Value = 883.0769
$ws.Range("I140").Value = 540
$ws.Range("K140").Value = 1620
$ws.Range("M140").Value = 3560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 20000.666
$ws.Range("I12").Value = 15001.5
$ws.Range("J12").Value = 29999
$ws.Range("K12").Value = 15001.5
$ws.Range("L12").Value = 29999
$ws.Range("M12").Value = -14861.5
$ws.Range("N12").Value = -30279

$ws.Range("H107").Value = 4831682
$ws.Range("I107").Value = 446.9091
$ws.Range("K107").Value = 446.9091
$ws.Range("M107").Value = 1473.0909

$ws.Range("H113").Value = 1243.3
$ws.Range("I113").Value = 1184.5
$ws.Range("J113").Value = 1331.5
$ws.Range("K113").Value = 1184.5
$ws.Range("L113").Value = 1331.5
$ws.Range("M113").Value = 985.5
$ws.Range("N113").Value = -5671.5

$ws.Range("H132").Value = 3898.5
$ws.Range("I132").Value = 2847.6365
$ws.Range("K132").Value = 8542.9095
$ws.Range("M132").Value = -6012.9095

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 63889.855
$ws.Range("J46").Value = 63889.855
$ws.Range("L46").Value = 63889.855
$ws.Range("N46").Value = -64351.855

$ws.Range("H134").Value = 63889.855
$ws.Range("J134").Value = 63889.855
$ws.Range("L134").Value = 191669.565
$ws.Range("N134").Value = -196739.565

$ws.Range("H136").Value = 5325.8335
$ws.Range("I136").Value = 1498.5
$ws.Range("K136").Value = 4495.5
$ws.Range("M136").Value = -1945.5
